# Auto-generated script applying the Ifrit_Profits (FFXIV leve-profit) market-data refresh
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 636.6667
$ws.Range("I53").Value = 336.2353
$ws.Range("J53").Value = 1147.4
$ws.Range("K53").Value = 336.2353
$ws.Range("L53").Value = 1147.4
$ws.Range("M53").Value = 300.7647
$ws.Range("N53").Value = -2421.4
$ws.Range("H74").Value = 4080
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 4400
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 4400
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -6272
$ws.Range("H76").Value = 3111.1
$ws.Range("I76").Value = 3111.1
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3111.1
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2796.1
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 4080
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 4400
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -31360
$ws.Range("H79").Value = 3111.1
$ws.Range("I79").Value = 3111.1
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3111.1
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2019.1
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 3599.4
$ws.Range("I116").Value = 3499.25
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 3499.25
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -57.25
$ws.Range("N116").Value = -10884

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5359.515
$ws.Range("I32").Value = 5355.125
$ws.Range("K32").Value = 5355.125
$ws.Range("M32").Value = -5068.125
$ws.Range("H74").Value = 2901102.5
$ws.Range("I74").Value = 3704153.8
$ws.Range("K74").Value = 3704153.8
$ws.Range("M74").Value = -3703279.8
$ws.Range("H77").Value = 2901102.5
$ws.Range("I77").Value = 3704153.8
$ws.Range("K77").Value = 18520769
$ws.Range("M77").Value = -18516401
$ws.Range("H97").Value = 655.3333
$ws.Range("I97").Value = 487.85715
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 487.85715
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = 8.14285000000001
$ws.Range("N97").Value = -3992

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 809.9
$ws.Range("I94").Value = 699.75
$ws.Range("J94").Value = 883.3333
$ws.Range("K94").Value = 699.75
$ws.Range("L94").Value = 883.3333
$ws.Range("M94").Value = -248.75
$ws.Range("N94").Value = -1785.3333
$ws.Range("H134").Value = 67803
$ws.Range("I134").Value = 81089.734
$ws.Range("J134").Value = 1369.3334
$ws.Range("K134").Value = 243269.202
$ws.Range("L134").Value = 4108.0002
$ws.Range("M134").Value = -240734.202
$ws.Range("N134").Value = -9178.0002

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2066
$ws.Range("I31").Value = 1272.5385
$ws.Range("J31").Value = 2639.0557
$ws.Range("K31").Value = 1272.5385
$ws.Range("L31").Value = 2639.0557
$ws.Range("M31").Value = -977.5385000000001
$ws.Range("N31").Value = -3229.0557
$ws.Range("H34").Value = 2066
$ws.Range("I34").Value = 1272.5385
$ws.Range("J34").Value = 2639.0557
$ws.Range("K34").Value = 1272.5385
$ws.Range("L34").Value = 2639.0557
$ws.Range("M34").Value = -1070.5385
$ws.Range("N34").Value = -3043.0557
$ws.Range("H41").Value = 13933.333
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 18400
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 18400
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -19256
$ws.Range("H134").Value = 1633.4193
$ws.Range("I134").Value = 1692.6428
$ws.Range("K134").Value = 5077.928400000001
$ws.Range("M134").Value = -2542.928400000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2636.3333
$ws.Range("I109").Value = 1427
$ws.Range("J109").Value = 2787.5
$ws.Range("K109").Value = 4281
$ws.Range("L109").Value = 8362.5
$ws.Range("M109").Value = -3241
$ws.Range("N109").Value = -10442.5
$ws.Range("H113").Value = 552.05884
$ws.Range("I113").Value = 540.1111
$ws.Range("J113").Value = 565.5
$ws.Range("K113").Value = 1620.3333
$ws.Range("L113").Value = 1696.5
$ws.Range("M113").Value = 549.6667000000002
$ws.Range("N113").Value = -6036.5
$ws.Range("H122").Value = 21009892
$ws.Range("I122").Value = 27778588
$ws.Range("K122").Value = 250007292
$ws.Range("M122").Value = -250004842
$ws.Range("H131").Value = 2532.2456
$ws.Range("I131").Value = 5893.3335
$ws.Range("J131").Value = 1902.0416
$ws.Range("K131").Value = 17680.0005
$ws.Range("L131").Value = 5706.1248
$ws.Range("M131").Value = -12640.0005
$ws.Range("N131").Value = -15786.1248

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 277501.5
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 550003
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 550003
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -551999
$ws.Range("H83").Value = 277501.5
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 550003
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 2750015
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -2759999
$ws.Range("H122").Value = 3786.8538
$ws.Range("I122").Value = 2926.16
$ws.Range("J122").Value = 5131.6875
$ws.Range("K122").Value = 8778.48
$ws.Range("L122").Value = 15395.0625
$ws.Range("M122").Value = -6328.48
$ws.Range("N122").Value = -20295.0625

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1394
$ws.Range("I40").Value = 1401.2307
$ws.Range("J40").Value = 1300
$ws.Range("K40").Value = 1401.2307
$ws.Range("L40").Value = 1300
$ws.Range("M40").Value = -1265.2307
$ws.Range("N40").Value = -1572
$ws.Range("H82").Value = 1826.8182
$ws.Range("I82").Value = 1282
$ws.Range("J82").Value = 2480.6
$ws.Range("K82").Value = 1282
$ws.Range("L82").Value = 2480.6
$ws.Range("M82").Value = -921
$ws.Range("N82").Value = -3202.6
$ws.Range("H85").Value = 1826.8182
$ws.Range("I85").Value = 1282
$ws.Range("J85").Value = 2480.6
$ws.Range("K85").Value = 1282
$ws.Range("L85").Value = 2480.6
$ws.Range("M85").Value = -34
$ws.Range("N85").Value = -4976.6

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10025.2
$ws.Range("J45").Value = 10025.2
$ws.Range("L45").Value = 10025.2
$ws.Range("N45").Value = -11007.2
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680
$ws.Range("H122").Value = 2034.8334
$ws.Range("I122").Value = 1175.8695
$ws.Range("J122").Value = 4857.143
$ws.Range("K122").Value = 3527.6085
$ws.Range("L122").Value = 14571.429
$ws.Range("M122").Value = -1077.6085
$ws.Range("N122").Value = -19471.429
$ws.Range("H132").Value = 9064.15
$ws.Range("I132").Value = 9515.833000000001
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 28547.499
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -26017.499
$ws.Range("N132").Value = -20057
